$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 14
$ws.Range("A4").Value = 24
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 3
$ws.Range("A7").Value = 2
$ws.Range("A8").Value = 1

$ws.Range("B2").Value = 15
$ws.Range("B3").Value = 15
$ws.Range("B4").Value = 12

$ws.Range("B4").Select()
